$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in the "Implementing routes" row description
$ws.Range("C36").Value = "All possible routes without authentication and authorization!"

# New block: row 42 - User Registration on Home Page (item 5)
$ws.Range("A42").Value = 5
$ws.Range("B42").Value = "Tempalte Driven Form"
$ws.Range("C42").Value = "User Registration on Home Page"

$ws.Range("C43").Value = "Service for using all REST endpoints"

$ws.Range("D44").Value = "json-server"
$ws.Range("C44").Value = "Create Fake REST end points"

$ws.Range("C45").Value = "Add New Shop Template Driven Form"

$ws.Range("C46").Value = "Service for add new shop"

# Update selection/view to match the new extent of data
$ws.Range("C47").Select()
